$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Protect column D (and any other cell) from Excel auto-converting
# numeric-looking text (e.g. "1.00", "0.0000133") into real numbers by
# forcing a Text number format before the writes, then clearing the
# formatting afterwards so the cell style stays untouched (matches the
# original workbook, which stores all of these as plain inline strings).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "70.279.16"
$ws.Range("E2").Value = "  +0.63%  "

$ws.Range("D3").Value = "3.760.53"
$ws.Range("E3").Value = "  -0.64%  "

$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "619.22"
$ws.Range("E5").Value = "  +0.55%  "

$ws.Range("D6").Value = "181.74"
$ws.Range("E6").Value = "  +2.60%  "

$ws.Range("D7").Value = "3.759.91"
$ws.Range("E7").Value = "  -0.63%  "

$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("E9").Value = "  -1.72%  "

$ws.Range("E10").Value = "  -1.03%  "

$ws.Range("D11").Value = "6.41"
$ws.Range("E11").Value = "  +0.47%  "

$ws.Range("E12").Value = "  -3.48%  "

$ws.Range("D13").Value = "40.24"

$ws.Range("E14").Value = "  -0.40%  "

$ws.Range("D15").Value = "4.383.21"
$ws.Range("E15").Value = "  -0.73%  "

$ws.Range("D16").Value = "3.759.44"
$ws.Range("E16").Value = "  -0.65%  "

$ws.Range("D17").Value = "70.321.89"
$ws.Range("E17").Value = "  +0.47%  "

$ws.Range("E18").Value = "  -2.10%  "

$ws.Range("D19").Value = "7.57"
$ws.Range("E19").Value = "  -0.28%  "

$ws.Range("D20").Value = "16.60"
$ws.Range("E20").Value = "  -0.62%  "

$ws.Range("D21").Value = "505.60"
$ws.Range("E21").Value = "  -2.48%  "

$ws.Range("D22").Value = "9.24"
$ws.Range("E22").Value = "  -1.64%  "

$ws.Range("D23").Value = "0.725"
$ws.Range("E23").Value = "  -2.07%  "

$ws.Range("D24").Value = "2.61"
$ws.Range("E24").Value = "  +4.58%  "

$ws.Range("D25").Value = "86.52"
$ws.Range("E25").Value = "  -2.67%  "

$ws.Range("D26").Value = "13.10"
$ws.Range("E26").Value = "  -3.53%  "

$ws.Range("D27").Value = "11.30"
$ws.Range("E27").Value = "  +3.72%  "

$ws.Range("D28").Value = "0.0000133"
$ws.Range("E28").Value = "  +4.23%  "

$ws.Range("E29").Value = "  +0.23%  "

$ws.Range("D30").Value = "2.51"
$ws.Range("E30").Value = "  -0.24%  "

$ws.Range("E31").Value = "  +2.33%  "

$ws.Range("D32").Value = "7.96"
$ws.Range("E32").Value = "  +1.05%  "

$ws.Range("D33").Value = "30.72"
$ws.Range("E33").Value = "  -4.41%  "

$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.01%  "

$ws.Range("E36").Value = "  +0.42%  "

$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("D38").Value = "0.352"
$ws.Range("E38").Value = "  +2.69%  "

$ws.Range("E39").Value = "  +5.70%  "

$ws.Range("D40").Value = "3.08"
$ws.Range("E40").Value = "  +13.88%  "

$ws.Range("E41").Value = "  -4.90%  "

$ws.Range("E42").Value = "  -2.91%  "

$ws.Range("D43").Value = "45.61"
$ws.Range("E43").Value = "  +2.10%  "

$ws.Range("D44").Value = "435.51"
$ws.Range("E44").Value = "  +2.07%  "

$ws.Range("D45").Value = "8.65"
$ws.Range("E45").Value = "  -2.13%  "

$ws.Range("D46").Value = "2.971.24"
$ws.Range("E46").Value = "  -4.59%  "

$ws.Range("E47").Value = "  -0.45%  "

$ws.Range("D48").Value = "27.51"
$ws.Range("E48").Value = "  -1.07%  "

$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "138.58"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("B50").Value = "USDe"
$ws.Range("C50").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D50").Value = "1.00"
$ws.Range("E50").Value = "  -0.03%  "

$ws.Range("D51").Value = "2.50"
$ws.Range("E51").Value = "  +0.71%  "

# Restore original (default) formatting on the cells we touched.
$ws.Range("D2:E51").ClearFormats()
